# DateFormatTests.xlsx edit
#
# Adds four new "Time" test rows (43-46) to the "Tests" sheet that exercise
# the "d \d\a\y\s h" / d "days" h (with and without am/pm) custom TEXT()
# format codes, plus the small cosmetic follow-up tweaks that came along
# with the author's commit (row heights on the "Flags" sheet, the "Tests"
# sheet's Value column width/selection, and the dimension growing to
# accommodate the new rows).

$wb = $excel.ActiveWorkbook
$flags = $wb.Worksheets.Item("Flags")
$tests = $wb.Worksheets.Item("Tests")

# ---------------------------------------------------------------------
# New data rows on the "Tests" sheet
# ---------------------------------------------------------------------
# Column layout: A = Result (expected) [=TEXT(C,B)], B = Format, C = Value,
# D = Categories. C43:C46 all reuse the same serial date/time value that
# rows 2 and 42 already use (1952-10-11 14:35:27 under the 1904 date
# system), and D43:D46 are categorized "Time" like row 42.

$rows = @(
    @{ Row = 43; Format = "d \d\a\y\s h" },
    @{ Row = 44; Format = "d ""days"" h" },
    @{ Row = 45; Format = "d \d\a\y\s h a/p" },
    @{ Row = 46; Format = "d ""days"" h am/pm" }
)

foreach ($item in $rows) {
    $r = $item.Row

    # Format column: plain text, "@" number format, new "Lucida Sans
    # Regular" font. Resetting to the Normal style first keeps the xf
    # minimal (no inherited alignment from the column's default style).
    $cellB = $tests.Range("B$r")
    $cellB.Style = "Normal"
    $cellB.NumberFormat = "@"
    $cellB.Font.Name = "Lucida Sans Regular"
    $cellB.Value = $item.Format

    # Value column: same serial date/time as the other "Time" rows, with
    # the shared custom date/time display format.
    $cellC = $tests.Range("C$r")
    $cellC.Value = 17816.607951388887
    $cellC.NumberFormat = "dd\-mmm\-yyyy\ hh:mm:ss.000"

    # Categories column.
    $tests.Range("D$r").Value = "Time"

    # Result column: TEXT() formula referencing this row's Format/Value.
    $tests.Range("A$r").Formula = "=TEXT(C$r,B$r)"
}

# Grow the used range / selection to include the new rows.
$tests.Activate()
$tests.Range("C46").Select()

# The Value column (C) was manually widened and is no longer "best fit".
$tests.Columns.Item(3).ColumnWidth = 48.8

# ---------------------------------------------------------------------
# "Flags" sheet row-height tweaks
# ---------------------------------------------------------------------
$flags.Rows.Item(1).RowHeight = 14
$flags.Rows.Item(2).RowHeight = 14
$flags.Rows.Item(3).RowHeight = 14
$flags.Rows.Item(4).RowHeight = 28
